$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 242, shifting existing rows 242-270 down to 243-271.
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new data record.
$ws.Range("A242").Value2 = 10
$ws.Range("B242").Value2 = "Vega Modelo de Temuco"
$ws.Range("C242").Value2 = "La Araucanía"
$ws.Range("D242").Value2 = 44505
$ws.Range("E242").Value2 = 9
$ws.Range("F242").Value2 = 100112024
$ws.Range("G242").Value2 = "Choclo"
$ws.Range("H242").Value2 = "Dulce o Americano"
$ws.Range("I242").Value2 = "Primera"
$ws.Range("J242").Value2 = 75
$ws.Range("K242").Value2 = 37000
$ws.Range("L242").Value2 = 38000
$ws.Range("M242").Value2 = 37467
$ws.Range("N242").Value2 = "$/malla 70 unidades"
$ws.Range("O242").Value2 = "Región de Arica y Parinacota"
$ws.Range("P242").Value2 = 535
$ws.Range("Q242").Value2 = 70
$ws.Range("R242").Value2 = "Hortaliza"
